$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unify header alignment first (G1:H1 were left-aligned) so the whole header
# row starts from one common style before layering on bold + borders.
$ws.Range("G1:H1").HorizontalAlignment = -4108

# Apply border + bold to the (now uniform) header row in one pass.
$headerRange = $ws.Range("A1:H1")
$headerRange.Borders.ColorIndex = 1
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true

# Insert the new RunNo column at the front; everything else shifts right.
$ws.Columns("A:A").Insert()

$ws.Range("A1").Value2 = "RunNo"
$ws.Range("A2").Value2 = 1
$ws.Range("A3").Value2 = 2

Write-Host "done"
